# Auto-generated script to apply numeric cell updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: H12, I12, J12, K12, L12, M12, N12
$ws.Range("H12").Value = 235.8
$ws.Range("I12").Value = 320
$ws.Range("J12").Value = 214.75
$ws.Range("K12").Value = 320
$ws.Range("L12").Value = 214.75
$ws.Range("M12").Value = -150
$ws.Range("N12").Value = -554.75
# Row 43: H43, J43, L43, N43
$ws.Range("H43").Value = 3454.5
$ws.Range("J43").Value = 3454.5
$ws.Range("L43").Value = 3454.5
$ws.Range("N43").Value = -3592.5
# Row 55: H55, I55, K55, M55
$ws.Range("H55").Value = 189.35
$ws.Range("I55").Value = 162.6
$ws.Range("K55").Value = 162.6
$ws.Range("M55").Value = 51.40000000000001
# Row 86: H86, I86, J86, K86, L86, M86, N86
$ws.Range("H86").Value = 1493.8
$ws.Range("I86").Value = 1348.8
$ws.Range("J86").Value = 1638.8
$ws.Range("K86").Value = 1348.8
$ws.Range("L86").Value = 1638.8
$ws.Range("M86").Value = -225.8
$ws.Range("N86").Value = -3884.8
# Row 87: H87, I87, K87, M87
$ws.Range("H87").Value = 12299.5
$ws.Range("I87").Value = 12299.5
$ws.Range("K87").Value = 12299.5
$ws.Range("M87").Value = -11051.5
# Row 89: H89, I89, J89, K89, L89, M89, N89
$ws.Range("H89").Value = 1493.8
$ws.Range("I89").Value = 1348.8
$ws.Range("J89").Value = 1638.8
$ws.Range("K89").Value = 6744
$ws.Range("L89").Value = 8194
$ws.Range("M89").Value = -1128
$ws.Range("N89").Value = -19426
# Row 90: H90, I90, K90, M90
$ws.Range("H90").Value = 12299.5
$ws.Range("I90").Value = 12299.5
$ws.Range("K90").Value = 36898.5
$ws.Range("M90").Value = -30658.5

$ws = $wb.Worksheets.Item("ARM")
# Row 55: H55, I55, J55, K55, L55, M55, N55
$ws.Range("H55").Value = 23186.428
$ws.Range("I55").Value = 20000
$ws.Range("J55").Value = 23717.5
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 23717.5
$ws.Range("M55").Value = -19685
$ws.Range("N55").Value = -24347.5
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 2291.9473
$ws.Range("I122").Value = 2066.4546
$ws.Range("K122").Value = 6199.3638
$ws.Range("M122").Value = -3749.3638
# Row 125: H125, J125, L125, N125
$ws.Range("H125").Value = 66916.86
$ws.Range("J125").Value = 66916.86
$ws.Range("L125").Value = 66916.86
$ws.Range("N125").Value = -76756.86

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31, I31, K31, M31
$ws.Range("H31").Value = 2279.6667
$ws.Range("I31").Value = 2515.6
$ws.Range("K31").Value = 2515.6
$ws.Range("M31").Value = -2220.6
# Row 34: H34, I34, K34, M34
$ws.Range("H34").Value = 2279.6667
$ws.Range("I34").Value = 2515.6
$ws.Range("K34").Value = 2515.6
$ws.Range("M34").Value = -2313.6
# Row 58: H58, I58, K58, M58
$ws.Range("H58").Value = 1994.2413
$ws.Range("I58").Value = 1871.5927
$ws.Range("K58").Value = 1871.5927
$ws.Range("M58").Value = -1668.5927
# Row 86: H86, J86, L86, N86
$ws.Range("H86").Value = 2993
$ws.Range("J86").Value = 3400
$ws.Range("L86").Value = 3400
$ws.Range("N86").Value = -5646
# Row 89: H89, J89, L89, N89
$ws.Range("H89").Value = 2993
$ws.Range("J89").Value = 3400
$ws.Range("L89").Value = 17000
$ws.Range("N89").Value = -28232
# Row 99: H99, I99, J99, K99, L99, M99, N99
$ws.Range("H99").Value = 2271
$ws.Range("I99").Value = 1895.1538
$ws.Range("J99").Value = 2813.889
$ws.Range("K99").Value = 1895.1538
$ws.Range("L99").Value = 2813.889
$ws.Range("M99").Value = -397.1538
$ws.Range("N99").Value = -5809.889
# Row 105: H105, I105, J105, K105, L105, M105, N105
$ws.Range("H105").Value = 1730
$ws.Range("I105").Value = 1254.125
$ws.Range("J105").Value = 2999
$ws.Range("K105").Value = 1254.125
$ws.Range("L105").Value = 2999
$ws.Range("M105").Value = 492.875
$ws.Range("N105").Value = -6493
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 2199.75
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250
# Row 126: H126, I126, J126, K126, L126, M126, N126
$ws.Range("H126").Value = 2271
$ws.Range("I126").Value = 1895.1538
$ws.Range("J126").Value = 2813.889
$ws.Range("K126").Value = 5685.4614
$ws.Range("L126").Value = 8441.667000000001
$ws.Range("M126").Value = -3215.4614
$ws.Range("N126").Value = -13381.667
# Row 132: H132, J132, L132, N132
$ws.Range("H132").Value = 5466.273
$ws.Range("J132").Value = 4995.091
$ws.Range("L132").Value = 14985.273
$ws.Range("N132").Value = -20045.273
# Row 134: H134, J134, L134, N134
$ws.Range("H134").Value = 3228183.2
$ws.Range("J134").Value = 12504630
$ws.Range("L134").Value = 37513890
$ws.Range("N134").Value = -37518960
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 1994.2413
$ws.Range("I136").Value = 1871.5927
$ws.Range("K136").Value = 5614.7781
$ws.Range("M136").Value = -3064.7781

$ws = $wb.Worksheets.Item("CUL")
# Row 5: H5, J5, L5, N5
$ws.Range("H5").Value = 716.5833
$ws.Range("J5").Value = 481.91666
$ws.Range("L5").Value = 1445.74998
$ws.Range("N5").Value = -1669.74998
# Row 33: H33, I33, J33, K33, L33, M33, N33
$ws.Range("H33").Value = 886
$ws.Range("I33").Value = 870
$ws.Range("J33").Value = 902
$ws.Range("K33").Value = 5220
$ws.Range("L33").Value = 5412
$ws.Range("M33").Value = -4937
$ws.Range("N33").Value = -5978
# Row 34: H34, J34, L34, N34
$ws.Range("H34").Value = 5828.4287
$ws.Range("J34").Value = 6199.846
$ws.Range("L34").Value = 18599.538
$ws.Range("N34").Value = -18767.538
# Row 44: H44, I44, K44, M44
$ws.Range("H44").Value = 861.4
$ws.Range("I44").Value = 825.75
$ws.Range("K44").Value = 2477.25
$ws.Range("M44").Value = -2079.25
# Row 46: H46, J46, L46, N46
$ws.Range("H46").Value = 94671.37
$ws.Range("J46").Value = 5645.5713
$ws.Range("L46").Value = 16936.7139
$ws.Range("N46").Value = -17118.7139
# Row 107: H107, J107, L107, N107
$ws.Range("H107").Value = 804.6
$ws.Range("J107").Value = 897.9167
$ws.Range("L107").Value = 2693.7501
$ws.Range("N107").Value = -6533.7501
# Row 135: H135, J135, L135, N135
$ws.Range("H135").Value = 716.5833
$ws.Range("J135").Value = 481.91666
$ws.Range("L135").Value = 4337.24994
$ws.Range("N135").Value = -9407.24994

$ws = $wb.Worksheets.Item("GSM")
# Row 126: H126, J126, L126, N126
$ws.Range("H126").Value = 5592.143
$ws.Range("J126").Value = 3479.4
$ws.Range("L126").Value = 10438.2
$ws.Range("N126").Value = -15378.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7: H7, I7, J7, K7, L7, M7, N7
$ws.Range("H7").Value = 3888.5557
$ws.Range("I7").Value = 3501.5
$ws.Range("J7").Value = 3999.1428
$ws.Range("K7").Value = 3501.5
$ws.Range("L7").Value = 3999.1428
$ws.Range("M7").Value = -3389.5
$ws.Range("N7").Value = -4223.1428
# Row 68: H68, I68, K68, M68
$ws.Range("H68").Value = 7333.6665
$ws.Range("I68").Value = 9499
$ws.Range("K68").Value = 9499
$ws.Range("M68").Value = -8750
# Row 71: H71, I71, K71, M71
$ws.Range("H71").Value = 7333.6665
$ws.Range("I71").Value = 9499
$ws.Range("K71").Value = 47495
$ws.Range("M71").Value = -43751
# Row 126: H126, I126, J126, K126, L126, M126, N126
$ws.Range("H126").Value = 3888.5557
$ws.Range("I126").Value = 3501.5
$ws.Range("J126").Value = 3999.1428
$ws.Range("K126").Value = 10504.5
$ws.Range("L126").Value = 11997.4284
$ws.Range("M126").Value = -8034.5
$ws.Range("N126").Value = -16937.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 62: H62, I62, K62, M62
$ws.Range("H62").Value = 4407
$ws.Range("I62").Value = 4860.5
$ws.Range("K62").Value = 4860.5
$ws.Range("M62").Value = -4236.5
# Row 65: H65, I65, K65, M65
$ws.Range("H65").Value = 4407
$ws.Range("I65").Value = 4860.5
$ws.Range("K65").Value = 24302.5
$ws.Range("M65").Value = -21182.5
# Row 100: H100, I100, K100, M100
$ws.Range("H100").Value = 1414.1818
$ws.Range("I100").Value = 1444.5
$ws.Range("K100").Value = 2889
$ws.Range("M100").Value = -2348
# Row 108: H108, J108, L108, N108
$ws.Range("H108").Value = 29850.666
$ws.Range("J108").Value = 29850.666
$ws.Range("L108").Value = 29850.666
$ws.Range("N108").Value = -37530.666
# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value = 2760.1428
$ws.Range("I126").Value = 2160.3333
$ws.Range("K126").Value = 6480.999899999999
$ws.Range("M126").Value = -4010.999899999999
